# Add a new "Project Link" slide (Title Only layout) at the end of the
# deck, give it a title + a textbox with the repo URL, and refresh the
# cached "datetimeFigureOut" placeholder text (6/12/2024 -> 6/13/2024)
# on the slide master and every slide layout, matching a same-day
# re-save by PowerPoint after the new slide was authored on 6/13/2024.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) New slide 13, "Title Only" layout (ppLayoutTitleOnly = 11), appended
#    after the current last slide (index 13 = Slides.Count + 1).
# ---------------------------------------------------------------------
$newIndex = $p.Slides.Count + 1
$slide = $p.Slides.Add($newIndex, 11)

# Title placeholder -> "Project Link"
$title = $slide.Shapes.Item(1)
$title.Name = "Title 1"
$title.TextFrame.TextRange.Text = "Project Link"
$title.TextFrame.TextRange.LanguageID = "en-US"

# Nudge the shape-id counter so the textbox we add next lands on id=4
# (matching the authored file, where shape id 3 was consumed elsewhere
# during editing) instead of the default id=3.
$spacer = $slide.Shapes.AddTextbox(1, 0, 0, 1, 1)
$spacer.Delete()

# Free-floating textbox with the project repository link.
$emuPerPt = 12700
$left   = 2314575 / $emuPerPt
$top    = 3059668 / $emuPerPt
$width  = 7954735 / $emuPerPt
$height = 369332  / $emuPerPt

$link = $slide.Shapes.AddTextbox(1, $left, $top, $width, $height)
$link.Name = "TextBox 3"
$link.Fill.Visible = 0
$link.TextFrame.WordWrap = -1
$link.TextFrame.AutoSize = 1
$link.TextFrame.TextRange.Text = "https://github.com/muralimura/APPSSDC-CS-Project1.git"
$link.TextFrame.TextRange.LanguageID = "en-IN"

# ---------------------------------------------------------------------
# 2) Refresh the cached date placeholder text on the slide master and
#    every slide layout (6/12/2024 -> 6/13/2024).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "6/12/2024") {
                $sh.TextFrame.TextRange.Text = "6/13/2024"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
